$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.740.66'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '1.601.33'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('E9').Value = '  +0.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.58'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0847'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.83%  '
$ws.Range('D12').Value = '1.826.30'
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').Value = '1.602.17'
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('E15').Value = '  +0.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.04'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('E17').Value = '  -1.17%  '
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.00%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '208.76'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('E21').Value = '  +0.60%  '
$ws.Range('E22').Value = '  -4.25%  '
$ws.Range('E23').Value = '  +0.96%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '143.63'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('E25').Value = '  +0.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.13'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.34'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.29%  '
$ws.Range('E29').Value = '  -1.39%  '
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('E31').Value = '  +0.82%  '
$ws.Range('E32').Value = '  +0.56%  '
$ws.Range('D33').Value = '1.280.74'
$ws.Range('E33').Value = '  -0.42%  '
$ws.Range('E34').Value = '  +1.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.22'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +15.77%  '
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.590'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.13%  '
$ws.Range('E38').Value = '  -1.06%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.47'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.61%  '
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.778'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '62.64'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.53%  '
$ws.Range('D44').Value = '1.738.55'
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '90.30'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.102'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.85%  '
$ws.Range('E48').Value = '  +0.74%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.52'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('E51').Value = '  +1.71%  '
